# "added number of families per region as of 2012"
#
# Renames the (previously unfilled) "Population" column to "Population (2012)",
# applies a thousands-separator number format to it, and adds a new
# "Number of Families (2012)" column populated with 2012 family-count data,
# plus a new (header-only) "Demographics" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("C1").Value = "Population (2012)"
$ws.Range("D1").Value = "Number of Families (2012)"
$ws.Range("E1").Value = "Demographics"

# Give the two new header cells (and the two trailing blank header cells
# the original author's formatting run also covered) the same centered
# "Heading 1" look already used by B1/C1.
$ws.Range("C1").Copy()
$ws.Range("D1:G1").PasteSpecial(-4122)

# --- Number of Families (2012) data -----------------------------------
$ws.Range("D2").Value = $null
$ws.Range("D3").Value = 375
$ws.Range("D4").Value = 1105
$ws.Range("D5").Value = 771
$ws.Range("D6").Value = 2386
$ws.Range("D7").Value = 3082
$ws.Range("D8").Value = 638
$ws.Range("D9").Value = 1165
$ws.Range("D10").Value = 1604
$ws.Range("D11").Value = $null
$ws.Range("D12").Value = 1577
$ws.Range("D13").Value = 902
$ws.Range("D14").Value = 772
$ws.Range("D15").Value = 976
$ws.Range("D16").Value = 1078
$ws.Range("D17").Value = 988
$ws.Range("D18").Value = 532
$ws.Range("D19").Value = 557

# --- Number formatting / alignment -------------------------------------
# Population (2012) column stays empty but gets a thousands-separator
# number format; Number of Families (2012) gets the same format, centered.
$ws.Range("C2:C19").NumberFormat = "#,##0"
$ws.Range("D2:D19").NumberFormat = "#,##0"
$ws.Range("D2:D19").HorizontalAlignment = -4108

# --- Column sizing -------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 36.75
$ws.Columns.Item(4).ColumnWidth = 33.25
$ws.Columns.Item(5).ColumnWidth = 18.917

# --- Page setup & selection ---------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("E2").Select() | Out-Null
